$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.157930970191956
$ws.Range("B1").Value = 2.148375272750854
$ws.Range("C1").Value = 10.29334354400635
$ws.Range("D1").Value = 2.549879550933838
$ws.Range("E1").Value = 1.266692876815796
